$p = $ppt.ActivePresentation

# --- 1) Update the cached "datetimeFigureOut" date field text from
#        9/25/2025 -> 9/27/2025 on the slide master and every slide layout.
$oldDate = "9/25/2025"
$newDate = "9/27/2025"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DateShape $p.SlideMaster.Shapes

# Every slide layout belonging to the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateShape $layouts.Item($L).Shapes
}

# --- 2) Reposition the "Figure : Distribution of Features in the Dataset"
#        caption textbox on slide 7.
$s7 = $p.Slides.Item(7)
$caption = $s7.Shapes.Item(2)
$caption.Left = 1960891
$caption.Top = 5734445
